$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '90.658.64'
$ws.Range('E2').Value = '  -0.66%  '

# Row 3
$ws.Range('D3').Value = '3.135.83'
$ws.Range('E3').Value = '  +0.55%  '

# Row 4
$ws.Range('E4').Value = '  -0.07%  '

# Row 5
$ws.Range('D5').Value = '''238.31'
$ws.Range('E5').Value = '  +8.82%  '

# Row 6
$ws.Range('D6').Value = '''631.34'
$ws.Range('E6').Value = '  +1.23%  '

# Row 7
$ws.Range('E7').Value = '  +10.56%  '

# Row 8
$ws.Range('D8').Value = '''0.356'
$ws.Range('E8').Value = '  -6.89%  '

# Row 9
$ws.Range('E9').Value = '  +0.01%  '

# Row 10
$ws.Range('D10').Value = '3.134.95'
$ws.Range('E10').Value = '  +0.56%  '

# Row 11
$ws.Range('D11').Value = '''0.724'
$ws.Range('E11').Value = '  +0.35%  '

# Row 12
$ws.Range('E12').Value = '  +4.33%  '

# Row 13
$ws.Range('D13').Value = '''36.69'
$ws.Range('E13').Value = '  +6.37%  '

# Row 14
$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').Value = '''0.0000243'
$ws.Range('E14').Value = '  -4.42%  '

# Row 15
$ws.Range('B15').Value = 'Toncoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D15').Value = '''5.62'
$ws.Range('E15').Value = '  +4.35%  '

# Row 16
$ws.Range('D16').Value = '90.510.73'
$ws.Range('E16').Value = '  -0.70%  '

# Row 17
$ws.Range('D17').Value = '3.709.14'
$ws.Range('E17').Value = '  +0.30%  '

# Row 18
$ws.Range('D18').Value = '3.122.60'
$ws.Range('E18').Value = '  +0.88%  '

# Row 19
$ws.Range('D19').Value = '''3.67'
$ws.Range('E19').Value = '  -2.56%  '

# Row 20
$ws.Range('D20').Value = '''14.44'
$ws.Range('E20').Value = '  +2.37%  '

# Row 21
$ws.Range('E21').Value = '  -3.61%  '

# Row 22
$ws.Range('D22').Value = '''453.89'
$ws.Range('E22').Value = '  +4.27%  '

# Row 23
$ws.Range('E23').Value = '  +10.41%  '

# Row 24
$ws.Range('D24').Value = '''9.10'
$ws.Range('E24').Value = '  +4.09%  '

# Row 25
$ws.Range('D25').Value = '''6.03'
$ws.Range('E25').Value = '  -1.81%  '

# Row 26
$ws.Range('E26').Value = '  +5.27%  '

# Row 27
$ws.Range('D27').Value = '''12.51'
$ws.Range('E27').Value = '  +2.87%  '

# Row 28
$ws.Range('E28').Value = '  -0.45%  '

# Row 29
$ws.Range('D29').Value = '''0.999'
$ws.Range('E29').Value = '  +0.09%  '

# Row 30
$ws.Range('D30').Value = '''9.88'
$ws.Range('E30').Value = '  +9.05%  '

# Row 31
$ws.Range('D31').Value = '''0.162'
$ws.Range('E31').Value = '  -3.26%  '

# Row 32
$ws.Range('D32').Value = '''27.59'
$ws.Range('E32').Value = '  +16.61%  '

# Row 33
$ws.Range('D33').Value = '''0.202'
$ws.Range('E33').Value = '  +34.31%  '

# Row 34
$ws.Range('D34').Value = '''3.87'
$ws.Range('E34').Value = '  +3.17%  '

# Row 35
$ws.Range('D35').Value = '''514.75'
$ws.Range('E35').Value = '  -2.46%  '

# Row 36
$ws.Range('E36').Value = '  +6.07%  '

# Row 37
$ws.Range('E37').Value = '  +1.13%  '

# Row 38
$ws.Range('D38').Value = '''1.94'
$ws.Range('E38').Value = '  +4.40%  '

# Row 39
$ws.Range('E39').Value = '  +3.85%  '

# Row 40
$ws.Range('D40').Value = '''0.798'
$ws.Range('E40').Value = '  -20.14%  '

# Row 41
$ws.Range('D41').Value = '''0.430'
$ws.Range('E41').Value = '  +13.79%  '

# Row 42
$ws.Range('D42').Value = '''0.0876'
$ws.Range('E42').Value = '  +4.43%  '

# Row 43
$ws.Range('E43').Value = '  -0.43%  '

# Row 44
$ws.Range('E44').Value = '  +0.02%  '

# Row 45
$ws.Range('D45').Value = '''3.37'
$ws.Range('E45').Value = '  +43.34%  '

# Row 46
$ws.Range('D46').Value = '''1.95'
$ws.Range('E46').Value = '  +2.99%  '

# Row 47
$ws.Range('D47').Value = '''0.706'
$ws.Range('E47').Value = '  +14.37%  '

# Row 48
$ws.Range('D48').Value = '''149.46'
$ws.Range('E48').Value = '  +1.37%  '

# Row 49
$ws.Range('B49').Value = 'Filecoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D49').Value = '''4.60'
$ws.Range('E49').Value = '  +10.19%  '

# Row 50
$ws.Range('B50').Value = 'OKB'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D50').Value = '''45.49'
$ws.Range('E50').Value = '  +3.37%  '

# Row 51
$ws.Range('D51').Value = '''1.35'
$ws.Range('E51').Value = '  +4.51%  '
